$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add row 26: 四方坪站 (shared string index 4)
$ws.Cells.Item(26, 1).Value = 45974
$ws.Cells.Item(26, 2).Value = "四方坪站"
$ws.Cells.Item(26, 3).Value = 8540.36
$ws.Cells.Item(26, 4).Value = 7632.91
$ws.Cells.Item(26, 5).Value = 2855.21
$ws.Cells.Item(26, 6).Value = 369

# Add row 27: 高岭站 (shared string index 5)
$ws.Cells.Item(27, 1).Value = 45974
$ws.Cells.Item(27, 2).Value = "高岭站"
$ws.Cells.Item(27, 3).Value = 4197.05
$ws.Cells.Item(27, 4).Value = 3633.19
$ws.Cells.Item(27, 5).Value = 1070.23
$ws.Cells.Item(27, 6).Value = 154

$ws.Range("H27").Select()
